$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.762.15"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "2.225.57"
$ws.Range("E3").Value = "  -5.18%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "294.08"
$ws.Range("E5").Value = "  -5.33%  "

$ws.Range("D6").Value = "84.51"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -3.40%  "

$ws.Range("D10").Value = "0.0789"
$ws.Range("E10").Value = "  -2.82%  "

$ws.Range("D11").Value = "29.87"
$ws.Range("E11").Value = "  -0.87%  "

$ws.Range("D12").Value = "47.88"
$ws.Range("E12").Value = "  -8.80%  "

$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.570.18"
$ws.Range("E14").Value = "  -5.10%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "6.32"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "14.11"
$ws.Range("E16").Value = "  -4.65%  "

$ws.Range("D17").Value = "2.223.56"
$ws.Range("E17").Value = "  -6.29%  "

$ws.Range("D18").Value = "0.720"
$ws.Range("E18").Value = "  -5.37%  "

$ws.Range("D19").Value = "39.693.92"
$ws.Range("E19").Value = "  -0.98%  "

$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").Value = "5.76"
$ws.Range("E21").Value = "  -5.45%  "

$ws.Range("D22").Value = "65.26"
$ws.Range("E22").Value = "  -4.32%  "

$ws.Range("D23").Value = "10.49"
$ws.Range("E23").Value = "  -1.67%  "

$ws.Range("D24").Value = "232.46"

$ws.Range("E26").Value = "  -5.70%  "

$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("D28").Value = "22.81"
$ws.Range("E28").Value = "  -3.67%  "

$ws.Range("E29").Value = "  -3.62%  "

$ws.Range("D30").Value = "9.19"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "153.56"
$ws.Range("E31").Value = "  -0.19%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "32.57"
$ws.Range("E32").Value = "  -6.33%  "

$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  -5.97%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.37"
$ws.Range("E35").Value = "  -5.31%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0703"
$ws.Range("E36").Value = "  -2.28%  "

$ws.Range("D37").Value = "16.04"
$ws.Range("E37").Value = "  +3.08%  "

$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("E39").Value = "  -1.16%  "

$ws.Range("D40").Value = "2.65"
$ws.Range("E40").Value = "  -5.35%  "

$ws.Range("D41").Value = "1.64"
$ws.Range("E41").Value = "  -4.71%  "

$ws.Range("E42").Value = "  -3.84%  "

$ws.Range("D43").Value = "1.947.04"
$ws.Range("E43").Value = "  -0.97%  "

$ws.Range("E44").Value = "  -3.56%  "

$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("D46").Value = "9.36"
$ws.Range("E46").Value = "  -1.13%  "

$ws.Range("D47").Value = "16.14"
$ws.Range("E47").Value = "  -8.12%  "

$ws.Range("D48").Value = "2.58"
$ws.Range("E48").Value = "  -4.60%  "

$ws.Range("D49").Value = "2.440.91"
$ws.Range("E49").Value = "  -4.94%  "

$ws.Range("D50").Value = "70.66"
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").Value = "88.94"
$ws.Range("E51").Value = "  -4.56%  "
